# Apply the "feat: add 2022-Q3 data" change:
#  1. Insert a new worksheet named "2022-Q3" right after "总计" (i.e. before "2022-Q2").
#  2. Populate it with the quarterly fund-holding detail rows.
#  3. Insert a corresponding new row at the top of the "总计" (summary) sheet,
#     shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet before the existing "2022-Q2" sheet.
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($existingQ2)
$newSheet.Name = "2022-Q3"

# Header row (row 1): B..H, column A stays empty on the header row.
$header = New-Object 'object[,]' 1,7
$header[0,0] = "基金代码"
$header[0,1] = "基金名称"
$header[0,2] = "基金规模"
$header[0,3] = "股票总仓位"
$header[0,4] = "仓位占比"
$header[0,5] = "持有市值(亿元)"
$header[0,6] = "仓位排名"
$newSheet.Range("B1:H1").Value = $header

# Data rows 2-8: make sure columns D:G (规模/仓位/占比/市值) are stored as TEXT,
# matching the source data (these come in as strings, not numbers) except the
# very last 持有市值 cell (G8) which is a genuine numeric 0.
$newSheet.Range("D2:G7").NumberFormat = "@"
$newSheet.Range("D8:F8").NumberFormat = "@"

$rows = New-Object 'object[,]' 7,8
# columns: A(idx) B(code) C(name) D(scale) E(position) F(ratio) G(marketvalue) H(rank)
$rows[0,0] = 0;  $rows[0,1] = "501011"; $rows[0,2] = "汇添富中证中药指数（LOF）A"; $rows[0,3] = "11.36"; $rows[0,4] = "94.73"; $rows[0,5] = "4.45"; $rows[0,6] = "0.5055"; $rows[0,7] = 4
$rows[1,0] = 1;  $rows[1,1] = "673110"; $rows[1,2] = "西部利得新润灵活配置混合A"; $rows[1,3] = "4.48";  $rows[1,4] = "76.88"; $rows[1,5] = "6.89"; $rows[1,6] = "0.3087"; $rows[1,7] = 2
$rows[2,0] = 2;  $rows[2,1] = "501012"; $rows[2,2] = "汇添富中证中药指数（LOF）C"; $rows[2,3] = "6.42";  $rows[2,4] = "94.73"; $rows[2,5] = "4.45"; $rows[2,6] = "0.2857"; $rows[2,7] = 4
$rows[3,0] = 3;  $rows[3,1] = "159647"; $rows[3,2] = "鹏华中证中药ETF";            $rows[3,3] = "6.16";  $rows[3,4] = "94.79"; $rows[3,5] = "4.50"; $rows[3,6] = "0.2772"; $rows[3,7] = 4
$rows[4,0] = 4;  $rows[4,1] = "562390"; $rows[4,2] = "银华中证中药ETF";            $rows[4,3] = "2.34";  $rows[4,4] = "98.09"; $rows[4,5] = "4.61"; $rows[4,6] = "0.1079"; $rows[4,7] = 4
$rows[5,0] = 5;  $rows[5,1] = "561510"; $rows[5,2] = "华泰柏瑞中证中药ETF";        $rows[5,3] = "2.02";  $rows[5,4] = "95.98"; $rows[5,5] = "4.54"; $rows[5,6] = "0.0917"; $rows[5,7] = 5
$rows[6,0] = 6;  $rows[6,1] = "015356"; $rows[6,2] = "西部利得新润灵活配置混合C"; $rows[6,3] = "0.00";  $rows[6,4] = "76.88"; $rows[6,5] = "6.89"; $rows[6,6] = 0;       $rows[6,7] = 2

$newSheet.Range("A2:H8").Value = $rows

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Insert a new row on the "总计" sheet for 2022-Q3, shifting the rest down.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A2:A2").EntireRow.Insert()

$summaryRow = New-Object 'object[,]' 1,4
$summaryRow[0,0] = 0
$summaryRow[0,1] = "2022-Q3"
$summaryRow[0,2] = 7
$summaryRow[0,3] = 1.58
$summary.Range("A2:D2").Value = $summaryRow

# Re-number the helper index column (A) for all the following rows (1..7) since
# the row insert shifts data but keeps old index values; they already happen to
# line up (0..6 for the 7 rows after the header), so nothing else to change.

$summary.Range("A1").Select()
